$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Rows 1-3: simple text replacements ---
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"

# --- Insert 10 new rows immediately after row 3 (before the old row 4) ---
# Each Rows.Add inserted before row 4 pushes the previously-added row down,
# so insert the values in reverse order to end up in the right sequence.
$newValues = @("105", "0.00002", "0.00006", "0.00004", "0.00001", "0.00003", "0.00004", "0.00006", "0.00343", "100.0")

for ($i = $newValues.Length - 1; $i -ge 0; $i--) {
    $beforeRow = $t.Rows.Item(4)
    $newRow = $t.Rows.Add($beforeRow)
    $newRow.Cells.Item(1).Range.Text = $newValues[$i]
}

# --- Collapse the three trailing multi-run rows into single-run rows ---
$rowCount = $t.Rows.Count
$t.Rows.Item($rowCount - 2).Cells.Item(1).Range.Text = "100"
$t.Rows.Item($rowCount - 1).Cells.Item(1).Range.Text = "0"
$t.Rows.Item($rowCount).Cells.Item(1).Range.Text = "103"

Write-Output "rows=$($t.Rows.Count)"
